$d = $word.ActiveDocument

# 1) "Historyjka użytkownika: ... sklepu Guitar Center" - merge split "Guitar" / " Center" runs
$d.Content.Find.Execute("sklepu Guitar Center", $true, $false, $false, $false, $false, $true, 1, $false, "sklepu Guitar Center", 2) | Out-Null

# 2) "Środowisko: Mozilla Firefox 72.0" - merge split "Firefox" run
$d.Content.Find.Execute("odowisko: Mozilla Firefox 72.0", $true, $false, $false, $false, $false, $true, 1, $false, "odowisko: Mozilla Firefox 72.0", 2) | Out-Null

# 3) "Każdy z rozpatrywanych ... w repozytorium github: " - merge split "Katalon" and "github" runs
$d.Content.Find.Execute("wykorzystane zostało narzędzie Katalon Studio do zautomatyzowania każdego przypadku testowego. Pliki skryptów zostały dołączone do plików projektu, ale również znajdują się w repozytorium github: ", $true, $false, $false, $false, $false, $true, 1, $false, "wykorzystane zostało narzędzie Katalon Studio do zautomatyzowania każdego przypadku testowego. Pliki skryptów zostały dołączone do plików projektu, ale również znajdują się w repozytorium github: ", 2) | Out-Null

# 3b) The merge above causes the two following hyperlink runs (which share identical
#     formatting with each other) to also get coalesced into a single run, which the
#     target does NOT want (they must stay as two separate "https://...GC.FR" + "."
#     runs). Force a harmless toggle on the final character so the engine re-splits
#     that trailing "." back into its own run, restoring the original structure.
$hyperlink = $d.Content
$hyperlink.Find.Execute("https://github.com/tomek270/GC.FR.") | Out-Null
$lastDot = $d.Range($hyperlink.End - 1, $hyperlink.End)
$lastDot.Font.Bold = 1
$lastDot.Font.Bold = 0

# 4) "Przetestowany został ... narzędzie Katalon Studio." - merge split "Katalon" run
$d.Content.Find.Execute("wykorzystane zostało narzędzie Katalon Studio.", $true, $false, $false, $false, $false, $true, 1, $false, "wykorzystane zostało narzędzie Katalon Studio.", 2) | Out-Null

# 5) "Testy zostały wykonane na platformie Mozilla Firefox 72.0" - merge split "Firefox" run
$d.Content.Find.Execute("na platformie Mozilla Firefox 72.0", $true, $false, $false, $false, $false, $true, 1, $false, "na platformie Mozilla Firefox 72.0", 2) | Out-Null

# 6) "Zrzut ekranu z narzędzia Katalon Studio przedstawiający ..." - merge split "Katalon" run
$d.Content.Find.Execute("z narzędzia Katalon Studio przedstawiający", $true, $false, $false, $false, $false, $true, 1, $false, "z narzędzia Katalon Studio przedstawiający", 2) | Out-Null

# 7) "Informacje o znalezionych bugach znajdują się w pliku bugi.xlsx" - merge split "bugach" run
$d.Content.Find.Execute("znalezionych bugach znajdują", $true, $false, $false, $false, $false, $true, 1, $false, "znalezionych bugach znajdują", 2) | Out-Null
